# Consolidation of Data Layer
#
# The shared-string table is being tidied up: "LESSON" and "ClassHOURS"
# were unused/duplicated concepts, so they are replaced in-place by the
# two new category labels "SUBJECTS" and "CALENDAR" (the only sheet cells
# that referenced those strings - C8 and C10 on Sheet1 - are repointed to
# the new text; Excel's own shared-string consolidation on save then drops
# the now-unreferenced "LESSON"/"ClassHOURS" entries from the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value  = "SUBJECTS"
$ws.Range("C10").Value = "CALENDAR"

# Reposition the view: scroll so row 4 is at the top and select C11
# (matches the sheetView/selection recorded in the saved workbook).
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
